# Week2TransitionData.xlsx edit: "tidying things for report, changed to
# backward difference"
#
# Updates the separation example table (rows 13-15, column C) from the old
# xa-at-separation values to the recomputed ("backward difference") values,
# and updates the transition/separation note in A17 with the new ReL
# threshold.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xa at separation values recalculated using the backward-difference method
$ws.Range("C13").Value = 0.26
$ws.Range("C14").Value = 0.26
$ws.Range("C15").Value = 0.26

# Updated transition note with the new ReL threshold
$ws.Range("A17").Value = "Transition occurs before separation at ReL >= 1.8E6"

# Move the active selection to A18, matching where the author left off editing
$ws.Range("A18").Select()
